$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @{}
$values = @{}

$names[2] = "ВердиоГаст® Растительный комплекс для улучшения пищеварения (БАД ),  капсулы"; $values[2] = 81632
$names[3] = "Чага (березовый гриб) 50г"; $values[3] = 15792
$names[4] = "Спорыш трава 50г"; $values[4] = 10789
$names[5] = "Солодка корни 50г"; $values[5] = 28830
$names[6] = "Пижма цветки 75г"; $values[6] = 14819
$names[7] = "Полынь горькая трава 50г"; $values[7] = 33866
$names[8] = "Мать-и-мачеха листья 35г"; $values[8] = 21800
$names[9] = "Чистотел трава 50г"; $values[9] = 13005
$names[10] = "Сенна листья 50г"; $values[10] = 16799
$names[11] = "Алтей корни 75г"; $values[11] = 4826
$names[12] = "Шиповник плоды низковитаминные 50г"; $values[12] = 24211
$names[13] = "Липа цветки 35г"; $values[13] = 16008
$names[14] = "Брусника листья 50г"; $values[14] = 15033
$names[15] = "Кукуруза столбики с рыльцами 40г"; $values[15] = 23942
$names[16] = "Дуба кора 75г"; $values[16] = 62559
$names[17] = "Мята перечная листья 50г"; $values[17] = 23025
$names[18] = "Эвкалипт прутовидный листья 75г"; $values[18] = 26353
$names[19] = "Сб. Грудной №4 50г"; $values[19] = 37026
$names[20] = "Ноготки цветки 50г"; $values[20] = 25955
$names[21] = "Ромашка цветки вн 50г"; $values[21] = 98343
$names[22] = "Багульник болотный побеги 50г"; $values[22] = 15083
$names[23] = "Сб. Фитонефрол (Урологический сбор) 50г"; $values[23] = 8991
$names[24] = "Укроп пахучий плоды 50г"; $values[24] = 69482
$names[25] = "Девясил корневища и корни 50г"; $values[25] = 20211
$names[26] = "Береза почки 50г"; $values[26] = 19859
$names[27] = "Чабрец трава 50г"; $values[27] = 23367
$names[28] = "Эрва шерстистая трава 30г"; $values[28] = 15230
$names[29] = "Валериана корневища с корнями 50г"; $values[29] = 22734
$names[30] = "Пустырник трава 50г"; $values[30] = 13384
$names[31] = "Боярышник плоды 75г"; $values[31] = 25218
$names[32] = "Бессмертник песчаный цветки 30г"; $values[32] = 32721
$names[33] = "Подорожник большой листья 50г"; $values[33] = 11116
$names[34] = "Шалфей листья 50г"; $values[34] = 43534
$names[35] = "Аир корневища 75г"; $values[35] = 9339
$names[36] = "Ламинарии слоевища (морская капуста) 100г"; $values[36] = 19918
$names[37] = "Лен семена 100г"; $values[37] = 75169
$names[38] = "Череда трава 50г"; $values[38] = 16019
$names[39] = "Крушина кора 50г"; $values[39] = 14300
$names[40] = "Рябина плоды 50г"; $values[40] = 2660
$names[41] = "Сб. Фитогепатол №2 (Желчегонный сбор №2) 35г"; $values[41] = 6511
$names[42] = "Толокнянка листья 50г"; $values[42] = 12207
$names[43] = "Можжевельник плоды 50г"; $values[43] = 20382
$names[44] = "Зверобой трава 50г"; $values[44] = 53513
$names[45] = "Тысячелистник трава 50г"; $values[45] = 23783
$names[46] = "Сб. Фитопектол №1 (Грудной сбор №1) 35г"; $values[46] = 9281
$names[47] = "Крапива листья 50г"; $values[47] = 25093
$names[48] = "Сб. Фитопектол №2 (Грудной сбор №2) 35г"; $values[48] = 12720
$names[49] = "Фп Фиточай `"Лактафитол`" (БАД) 20х1,5 г"; $values[49] = 13229
$names[50] = "Фп Детский травяной чай `"ФармаЦветик® для иммунитета`" 20х1,5 г"; $values[50] = 2538
$names[51] = "Фп Детский травяной чай `"ФармаЦветик®  при простуде`" 20х1,5 г"; $values[51] = 4331
$names[52] = "Фп Детский травяной чай `"ФармаЦветик® для животика`" 20х1,5 г"; $values[52] = 4610
$names[53] = "Фп Детский травяной чай `"ФармаЦветик® для спокойного сна`" 20х1,5 г"; $values[53] = 7408
$names[54] = "Фп `"ВердиоГаст® Фиточай для улучшения пищеварения с зеленым чаем`"(БАД) 20*1,5г"; $values[54] = 8250
$names[55] = "Фп `"ВердиоГаст® Фиточай для улучшения пищеварения с черным чаем`" (БАД) 20*1,5г"; $values[55] = 9720
$names[56] = "Фп Фиточай `"Баланс`" (БАД) 20х2,0 г"; $values[56] = 90
$names[57] = "Фп Фиточай `"Дивный вечер`" (БАД) 20х2,0 г"
$names[58] = "Фп `"Щедрость природы® Фиточай для иммунитета`" 20х2,0 г"; $values[58] = 1134
$names[59] = "Фп `"Щедрость природы® Фиточай кардиологический`" 20х2,0 г"; $values[59] = 1476
$names[60] = "Фп `"Щедрость природы® Фиточай при простуде`" 20х2,0 г"; $values[60] = 1098
$names[61] = "Фп `"Щедрость природы® Фиточай успокоительный`"20х2,0 г"; $values[61] = 2502
$names[62] = "Фп Шалфей листья 20х1,5г"; $values[62] = 112794
$names[63] = "Фп `"Щедрость природы® Фиточай диабетический`" 20х2,0 г"; $values[63] = 1062
$names[64] = "Фп Сб. Фитогепатол №3 (Желчегонный сбор №3) 20x2,0г"; $values[64] = 58559
$names[65] = "Фп Сб. Грудной №4 20x2,0г"; $values[65] = 566985
$names[66] = "Фп Сб. Фитоседан №2 (Успокоительный сбор №2) 20x2,0г"; $values[66] = 44649
$names[67] = "Фп Сб. Бруснивер 20x2,0г"; $values[67] = 172263
$names[68] = "Фп Сб. Проктофитол (Противогеморроидальный сбор) 20х2,0г"; $values[68] = 20620
$names[69] = "Фп Сб. Желудочный №3 20x2,0г"; $values[69] = 23043
$names[70] = "Фп Толокнянка листья 20x1,5г"; $values[70] = 32328
$names[71] = "Фп Липа цветки 20x1,5г"; $values[71] = 65032
$names[72] = "Фп Сб. Фитонефрол (Урологический сбор) 20x2,0г"; $values[72] = 171139
$names[73] = "Фп Сб. Фитогастрол (Желудочно-кишечный сбор) 20x2,0г"; $values[73] = 79425
$names[74] = "Фп Аир корневища 20x1,5г"; $values[74] = 3784
$names[75] = "Фп Мелисса лекарственная трава 20x1,5г"; $values[75] = 38989
$names[76] = "Фп Ромашка цветки 20x1,5г"; $values[76] = 1412417
$names[77] = "Фп Боярышник плоды 20х3,0г"; $values[77] = 18260
$names[78] = "Фп Череда трава 20х1,5г"; $values[78] = 47682
$names[79] = "Фп Сб. Элекасол 20x2,0г"; $values[79] = 43574
$names[80] = "Фп Сенна листья 20x1,5г"; $values[80] = 72654
$names[81] = "Фп Пижма цветки 20х1,5г"; $values[81] = 6402
$names[82] = "Фп Сб. Фитоседан №3 (Успокоительный сбор №3) 20х2,0г"; $values[82] = 92049
$names[83] = "Фп Шиповник плоды 20х2,0г"; $values[83] = 59609
$names[84] = "Фп Фиточай `"Тибетский`" (БАД) 20х2,0  г"; $values[84] = 9270
$names[85] = "Фп Фиточай `"Опалиховский`" (БАД) 20х2,0 г"; $values[85] = 5814
$names[86] = "Фп Зверобой трава 20x1,5г"; $values[86] = 64537
$names[87] = "Фп `"Щедрость природы® Фиточай очищающий`" 20х2,0 г"; $values[87] = 1818
$names[88] = "Фп Подорожник листья 20x1,5г"; $values[88] = 35425
$names[89] = "Фп Сб. Арфазетин-Э 20x2,0г"; $values[89] = 48311
$names[90] = "Фп Брусника листья 20х1,5г"; $values[90] = 91259
$names[91] = "Фп Пустырник трава 20x1,5г"; $values[91] = 51837
$names[92] = "Фп Мята перечная листья 20x1,5г"; $values[92] = 81928
$names[93] = "Фп Чистотел трава 20х1,5г"; $values[93] = 37878
$names[94] = "Фп `"Щедрость природы® Фиточай для пищеварения`" 20х2,0 г"; $values[94] = 1854
$names[95] = "Фп Чабрец трава 20x1,5 г"; $values[95] = 94227
$names[96] = "Фп Душица трава 20x1,5г"; $values[96] = 37962
$names[97] = "Фп Крапива листья 20x1,5г"; $values[97] = 89152
$names[98] = "Фп Хвощ полевой трава 20х1,5г"; $values[98] = 40786
$names[99] = "Фп Пастушья сумка трава 20х1,5г"; $values[99] = 7378
$names[100] = "Фп Береза листья 20x1,5г"; $values[100] = 7006
$names[101] = "Фп Золототысячник трава 20х1,5г"; $values[101] = 6901
$names[102] = "Фп Фиалка трехцветная трава 20x1,5г"; $values[102] = 6520
$names[103] = "Фп Ольха соплодия 20х1,5г"; $values[103] = 6523
$names[104] = "Фп Ноготки цветки 20x1,5г"; $values[104] = 19632
$names[105] = "Фп Кровохлебка корневища и корни 20x1,5г"; $values[105] = 4084
$names[106] = "Фп Почечный чай листья 20x1,5г"; $values[106] = 51288
$names[107] = "Фп Валериана корневища с корнями 20x1,5г"; $values[107] = 17979
$names[108] = "Фп Лапчатка корневища 20x2,5г"; $values[108] = 3718
$names[109] = "Фп Девясил корневища и корни 20х1,5г"; $values[109] = 15090
$names[110] = "Фп Тысячелистник трава 20x1,5г"; $values[110] = 22156
$names[111] = "Фп Крушина кора 20x1,5г"; $values[111] = 12969
$names[112] = "Фп Бадан корневища 20x1,5г"; $values[112] = 2869

for ($r = 2; $r -le 112; $r++) {
    $ws.Cells.Item($r, 1).Value = $names[$r]
    if ($values.ContainsKey($r)) {
        $ws.Cells.Item($r, 2).Value = $values[$r]
    } else {
        $ws.Cells.Item($r, 2).Value = ""
    }
}

$ws.Range("A13").Select()
